# edit.ps1 - apply the changes described by the commit:
#   "add new figures, delete old"
#
# Concretely (per the OOXML diff):
#   1. The "Date Placeholder" field cached on the Slide Master and on every
#      Slide Layout changes its displayed text from 1/18/22 to 2/10/22
#      (the deck was re-saved on a later day; PowerPoint re-stamps the
#      cached text of the datetimeFigureOut field on the footer/date
#      placeholders of the master + every layout).
#   2. On slide 1, the small "color (dwl)" caption (previously split across
#      three runs: "color (", "dwl", ")") is normalized into a single run
#      with text "color (dwl)".

$p = $ppt.ActivePresentation
$newDate = "2/10/22"

function Set-DatePlaceholderText($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# -- 1. Update the cached date text on the Slide Master ---------------------
Set-DatePlaceholderText $p.SlideMaster.Shapes

# -- ... and on every Slide Layout hanging off that master ------------------
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $lyt = $layouts.Item($i)
    Set-DatePlaceholderText $lyt.Shapes
}

# -- 2. Merge the "color (dwl)" runs on slide 1 ------------------------------
$slide = $p.Slides.Item(1)
for ($k = 1; $k -le $slide.Shapes.Count; $k++) {
    $shp = $slide.Shapes.Item($k)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf("color (")
        if ($idx -ge 0) {
            # TextRange character offsets are 1-based.
            $start = $idx + 1
            $len = "color (dwl)".Length
            $sub = $tr.Characters($start, $len)
            $sub.Text = "color (dwl)"
        }
    }
}

Write-Output "done"
